$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("H(f) deriv no comp")
$ws2 = $wb.Worksheets.Item("Z(f) no comp")

# --- Sheet1 "H(f) deriv no comp": add column F (angular-frequency x RC calc) ---
$ws1.Range("F2").Formula = "=2*PI()*A2*1000*C2/2/10^6"
$ws1.Range("F2").ClearFormats()

$ws1.Range("F3:F13").Formula = "=2*PI()*A3*1000*C3/2/10^6"
$ws1.Range("F3:F13").ClearFormats()

# --- Sheet2 "Z(f) no comp": flip sign of the phase column D (rows 2-11) ---
$ws2.Range("D2").Formula = "=-96"
$ws2.Range("D3").Value = -90
$ws2.Range("D4").Value = -90
$ws2.Range("D5").Value = -89
$ws2.Range("D6").Value = -89
$ws2.Range("D7").Value = -90
$ws2.Range("D8").Value = -89
$ws2.Range("D9").Value = -88
$ws2.Range("D10").Value = -90
$ws2.Range("D11").Value = -91

# --- Selections / active tab: sheet2 becomes the active tab, selections updated ---
$ws1.Range("C10").Select()
$ws2.Activate()
$ws2.Range("D2:D11").Select()
